# Recomputed the results with the updated IKIEF data
#
# Each model-coefficient sheet (one per IK/Q/P curve) stores its fitted
# parameters in column B, rows 2-10 (offset, offset_datum, slope,
# temp_mean, temp_delta, time_offset, method, temp_ref, model_std) plus
# a 'gewijzigd' (last recomputed) timestamp in B11. Re-running the fit
# against the refreshed IKIEF input data changes the coefficients and
# bumps every sheet's timestamp to the new run.
#
# Values are written as plain (non-exponential) decimal strings: the
# host PowerShell-style parser here doesn't accept `1e-06`-style numeric
# literals, and round-tripping very small magnitudes through Excel as an
# exponential *string* causes Range.Value's auto-detection to stamp a
# Scientific NumberFormat onto the cell -- which the source workbook
# does not have. Plain-decimal strings parse to the identical double
# and leave formatting untouched.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("IK106")
$ws.Range("B2").Value = "-0.01411569026034072"
$ws.Range("B4").Value = "-0.000001006331873506303"
$ws.Range("B5").Value = "11.90303794667552"
$ws.Range("B6").Value = "6.707037223142624"
$ws.Range("B7").Value = "160.4616897849417"
$ws.Range("B10").Value = "0.3153588384032504"
$ws.Range("B11").Value = "45659.66900156771"

$ws = $wb.Worksheets.Item("Q100")
$ws.Range("B2").Value = "-0.01123846348618225"
$ws.Range("B4").Value = "-0.0000003627530010202139"
$ws.Range("B6").Value = "5.192187464916277"
$ws.Range("B7").Value = "156.0150837976051"
$ws.Range("B10").Value = "0.2761139415964836"
$ws.Range("B11").Value = "45659.66747815972"

$ws = $wb.Worksheets.Item("Q200")
$ws.Range("B11").Value = "45659.66756796296"

$ws = $wb.Worksheets.Item("Q300")
$ws.Range("B11").Value = "45659.66764251157"

$ws = $wb.Worksheets.Item("Q400")
$ws.Range("B11").Value = "45659.66773114583"

$ws = $wb.Worksheets.Item("Q500")
$ws.Range("B11").Value = "45659.66781262732"

$ws = $wb.Worksheets.Item("Q600")
$ws.Range("B11").Value = "45659.66790210648"

$ws = $wb.Worksheets.Item("P100")
$ws.Range("B11").Value = "45659.66798706019"

$ws = $wb.Worksheets.Item("P200")
$ws.Range("B11").Value = "45659.66808515046"

$ws = $wb.Worksheets.Item("P300")
$ws.Range("B11").Value = "45659.66814591435"

$ws = $wb.Worksheets.Item("P400")
$ws.Range("B11").Value = "45659.66822160879"

$ws = $wb.Worksheets.Item("P500")
$ws.Range("B11").Value = "45659.66829737269"

$ws = $wb.Worksheets.Item("P600")
$ws.Range("B11").Value = "45659.66837484953"

$ws = $wb.Worksheets.Item("IK91")
$ws.Range("B2").Value = "-0.03457833426573426"
$ws.Range("B4").Value = "-0.000001174062977483533"
$ws.Range("B5").Value = "13.07242886253713"
$ws.Range("B6").Value = "3.948032089390761"
$ws.Range("B7").Value = "183.7642883191519"
$ws.Range("B10").Value = "0.4715875904658048"
$ws.Range("B11").Value = "45659.66842083333"

$ws = $wb.Worksheets.Item("IK92")
$ws.Range("B2").Value = "-0.01226351728587713"
$ws.Range("B4").Value = "-0.000000964532815643391"
$ws.Range("B5").Value = "11.47904678264748"
$ws.Range("B6").Value = "6.176698667876619"
$ws.Range("B7").Value = "154.6122372058187"
$ws.Range("B10").Value = "0.1501591720239783"
$ws.Range("B11").Value = "45659.66846332176"

$ws = $wb.Worksheets.Item("IK93")
$ws.Range("B2").Value = "-0.01421276197693501"
$ws.Range("B4").Value = "-0.0000003069798351770336"
$ws.Range("B5").Value = "12.98193427849214"
$ws.Range("B6").Value = "7.463455860761829"
$ws.Range("B7").Value = "156.4442095998867"
$ws.Range("B10").Value = "0.1778421123817109"
$ws.Range("B11").Value = "45659.6685094213"

$ws = $wb.Worksheets.Item("IK94")
$ws.Range("B2").Value = "-0.01763472614621083"
$ws.Range("B4").Value = "-0.0000004097105753199457"
$ws.Range("B5").Value = "12.17380490732406"
$ws.Range("B6").Value = "7.55327129542511"
$ws.Range("B7").Value = "159.8862662245484"
$ws.Range("B10").Value = "0.2797756722834719"
$ws.Range("B11").Value = "45659.66855912037"

$ws = $wb.Worksheets.Item("IK95")
$ws.Range("B2").Value = "-0.01898798764021069"
$ws.Range("B4").Value = "-0.000001808033037566575"
$ws.Range("B5").Value = "12.81055004721435"
$ws.Range("B6").Value = "6.127366898086563"
$ws.Range("B7").Value = "148.0958320010689"
$ws.Range("B10").Value = "0.6674640416494141"
$ws.Range("B11").Value = "45659.66860541666"

$ws = $wb.Worksheets.Item("IK96")
$ws.Range("B2").Value = "-0.03485481994966595"
$ws.Range("B4").Value = "-0.000001162993107433335"
$ws.Range("B5").Value = "13.67921483785185"
$ws.Range("B6").Value = "3.067044428786604"
$ws.Range("B7").Value = "198.5677800483085"
$ws.Range("B10").Value = "0.2958347205778908"
$ws.Range("B11").Value = "45659.66864973379"

$ws = $wb.Worksheets.Item("IK101")
$ws.Range("B2").Value = "-0.03070645667343086"
$ws.Range("B4").Value = "-0.0000008280779236216774"
$ws.Range("B5").Value = "11.86272102785905"
$ws.Range("B6").Value = "4.956903238369805"
$ws.Range("B7").Value = "163.3719017798673"
$ws.Range("B10").Value = "0.3269768806603147"
$ws.Range("B11").Value = "45659.66869983797"

$ws = $wb.Worksheets.Item("IK102")
$ws.Range("B2").Value = "-0.01333316970482867"
$ws.Range("B4").Value = "-0.0000000001"
$ws.Range("B5").Value = "12.63675001236583"
$ws.Range("B6").Value = "6.286235868673095"
$ws.Range("B7").Value = "152.5821117267674"
$ws.Range("B10").Value = "0.2025706647313184"
$ws.Range("B11").Value = "45659.6687655787"

$ws = $wb.Worksheets.Item("IK103")
$ws.Range("B2").Value = "-0.01163233323401294"
$ws.Range("B4").Value = "-0.0000005047733073406892"
$ws.Range("B5").Value = "12.10879048092163"
$ws.Range("B6").Value = "4.51822951852187"
$ws.Range("B7").Value = "149.7057642482395"
$ws.Range("B10").Value = "0.2358357878301115"
$ws.Range("B11").Value = "45659.66883487268"

$ws = $wb.Worksheets.Item("IK104")
$ws.Range("B2").Value = "-0.01410927240071667"
$ws.Range("B4").Value = "-0.0000003807398764706498"
$ws.Range("B5").Value = "12.0774177547363"
$ws.Range("B6").Value = "6.97467375009463"
$ws.Range("B7").Value = "159.2688855138107"
$ws.Range("B10").Value = "0.4966794181628044"
$ws.Range("B11").Value = "45659.66888287037"

$ws = $wb.Worksheets.Item("IK105")
$ws.Range("B2").Value = "-0.01423617941937147"
$ws.Range("B4").Value = "-0.0000000001000000000000001"
$ws.Range("B5").Value = "11.99260469711941"
$ws.Range("B6").Value = "7.271024211687462"
$ws.Range("B7").Value = "154.6154940025039"
$ws.Range("B10").Value = "0.2292007559872446"
$ws.Range("B11").Value = "45659.6689434838"
